$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New dt_insertion timestamp (column H) applied to every data row (2-21)
$newDate = 45491.00302083333

# Per-row updates for position (E), points (F), matches (G)
# Row => E, F, G (use $null to leave a value unchanged)
$updates = @{
    2  = @($null, 36, 17)
    3  = @($null, $null, 17)
    4  = @($null, $null, $null)
    5  = @(5, $null, $null)
    6  = @(4, 30, 17)
    7  = @(6, $null, $null)
    8  = @($null, 29, 16)
    9  = @($null, $null, $null)
    10 = @(10, $null, $null)
    11 = @(11, $null, $null)
    12 = @(9, 23, $null)
    13 = @($null, $null, $null)
    14 = @($null, $null, $null)
    15 = @($null, $null, $null)
    16 = @($null, $null, 17)
    17 = @($null, $null, $null)
    18 = @($null, $null, $null)
    19 = @(18, $null, 15)
    20 = @(19, 11, $null)
    21 = @($null, $null, $null)
}

foreach ($row in 2..21) {
    $vals = $updates[$row]

    if ($null -ne $vals[0]) {
        $ws.Cells.Item($row, 5).Value = $vals[0]
    }
    if ($null -ne $vals[1]) {
        $ws.Cells.Item($row, 6).Value = $vals[1]
    }
    if ($null -ne $vals[2]) {
        $ws.Cells.Item($row, 7).Value = $vals[2]
    }

    $ws.Cells.Item($row, 8).Value = $newDate
}
